# Sam Curran / Chennai Super Kings activity data refresh.
# The source table (rows 2-11) is reordered/updated and a new row (12) is
# appended, matching the "updated activity till excel form" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 12 needs player/team labels like the existing rows (trailing
# char is U+00A0, matching the rest of column A).
$ws.Range("A12").NumberFormat = "@"
$ws.Range("A12").Value = "Sam Curran" + [char]0x00A0
$ws.Range("B12").NumberFormat = "@"
$ws.Range("B12").Value = "Chennai Super Kings"

# runs / balls / fours / sixes for rows 2-12 (keep stored-as-text like the
# rest of the sheet).
$ws.Range("C2:F12").NumberFormat = "@"

$ws.Range("C2").Value = "52"
$ws.Range("D2").Value = "47"
$ws.Range("E2").Value = "4"
$ws.Range("F2").Value = "2"

$ws.Range("C3").Value = "13"
$ws.Range("D3").Value = "14"
$ws.Range("E3").Value = "1"
$ws.Range("F3").Value = "0"

$ws.Range("C4").Value = "17"
$ws.Range("D4").Value = "11"
$ws.Range("E4").Value = "1"
$ws.Range("F4").Value = "1"

$ws.Range("C5").Value = "17"
$ws.Range("D5").Value = "6"
$ws.Range("E5").Value = "1"
$ws.Range("F5").Value = "2"

$ws.Range("C6").Value = "22"
$ws.Range("D6").Value = "25"
$ws.Range("E6").Value = "1"
$ws.Range("F6").Value = "1"

$ws.Range("C7").Value = "31"
$ws.Range("D7").Value = "21"
$ws.Range("E7").Value = "3"
$ws.Range("F7").Value = "2"

$ws.Range("C8").Value = "0"
$ws.Range("D8").Value = "3"
$ws.Range("E8").Value = "0"
$ws.Range("F8").Value = "0"

$ws.Range("C9").Value = "15"
$ws.Range("D9").Value = "5"
$ws.Range("E9").Value = "0"
$ws.Range("F9").Value = "2"

$ws.Range("C10").Value = "18"
$ws.Range("D10").Value = "6"
$ws.Range("E10").Value = "1"
$ws.Range("F10").Value = "2"

$ws.Range("C11").Value = "0"
$ws.Range("D11").Value = "1"
$ws.Range("E11").Value = "0"
$ws.Range("F11").Value = "0"

$ws.Range("C12").Value = "1"
$ws.Range("D12").Value = "2"
$ws.Range("E12").Value = "0"
$ws.Range("F12").Value = "0"
